# Apply the crypto-price refresh described by the commit diff.
# Price-looking D-column values are forced to Text via NumberFormat "@"
# so Excel does not silently coerce them into numbers (the source sheet
# stores every Price/Volume cell as text, including plain-looking ones).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.849.92"
$ws.Range("E2").Value = "  +0.19%  "
# Row 3
$ws.Range("D3").Value = "1.976.31"
$ws.Range("E3").Value = "  +0.44%  "
# Row 4
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.19"
$ws.Range("E5").Value = "  +0.12%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  +1.30%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.86"
$ws.Range("E7").Value = "  +2.61%  "
# Row 8
$ws.Range("E8").Value = "  +0.01%  "
# Row 9
$ws.Range("E9").Value = "  +1.86%  "
# Row 10
$ws.Range("E10").Value = "  -1.48%  "
# Row 11
$ws.Range("E11").Value = "  +0.68%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.62"
$ws.Range("E12").Value = "  +6.03%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.844"
$ws.Range("E13").Value = "  +1.79%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.98"
$ws.Range("E14").Value = "  -2.23%  "
# Row 15
$ws.Range("D15").Value = "2.269.66"
$ws.Range("E15").Value = "  +0.47%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.43"
$ws.Range("E16").Value = "  +2.76%  "
# Row 17
$ws.Range("D17").Value = "1.978.18"
$ws.Range("E17").Value = "  -0.01%  "
# Row 18
$ws.Range("D18").Value = "36.781.53"
$ws.Range("E18").Value = "  +0.24%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.12"
$ws.Range("E19").Value = "  +0.30%  "
# Row 20
$ws.Range("E20").Value = "  -0.17%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.15"
$ws.Range("E21").Value = "  +1.07%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.64"
$ws.Range("E22").Value = "  +0.03%  "
# Row 23
$ws.Range("E23").Value = "  +0.05%  "
# Row 24
$ws.Range("E24").Value = "  +1.44%  "
# Row 25
$ws.Range("E25").Value = "  +0.22%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.145"
$ws.Range("E26").Value = "  +1.66%  "
# Row 27
$ws.Range("E27").Value = "  -0.99%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.40"
$ws.Range("E28").Value = "  +1.74%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.45"
$ws.Range("E29").Value = "  +0.07%  "
# Row 30
$ws.Range("E30").Value = "  +18.23%  "
# Row 31
$ws.Range("E31").Value = "  +1.43%  "
# Row 32
$ws.Range("E32").Value = "  +2.40%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0619"
$ws.Range("E33").Value = "  -0.23%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.51"
$ws.Range("E34").Value = "  +4.87%  "
# Row 35
$ws.Range("E35").Value = "  +0.15%  "
# Row 36
$ws.Range("E36").Value = "  +0.01%  "
# Row 37
$ws.Range("E37").Value = "  -2.28%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.48"
$ws.Range("E39").Value = "  -10.04%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0976"
$ws.Range("E40").Value = "  -2.86%  "
# Row 41
$ws.Range("E41").Value = "  +0.84%  "
# Row 42
$ws.Range("E42").Value = "  +0.27%  "
# Row 43
$ws.Range("E43").Value = "  -0.30%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.11"
$ws.Range("E44").Value = "  -0.14%  "
# Row 45
$ws.Range("D45").Value = "1.367.24"
$ws.Range("E45").Value = "  +0.08%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.73"
$ws.Range("E46").Value = "  +1.82%  "
# Row 47
$ws.Range("E47").Value = "  -0.61%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.23"
$ws.Range("E48").Value = "  +0.82%  "
# Row 49
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.81"
$ws.Range("E49").Value = "  -0.95%  "
# Row 50
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.28"
$ws.Range("E50").Value = "  +5.48%  "
# Row 51
$ws.Range("E51").Value = "  +8.28%  "
